$wb = $excel.ActiveWorkbook

$wsPresupuesto = $wb.Worksheets.Item("Presupuesto")
$wsCompras = $wb.Worksheets.Item("Compras")
$wsInventario = $wb.Worksheets.Item("Inventario_cero_coste")

# --- Inventario_cero_coste: new inventory row for the purchased servo (row 12) ---
# Add the product link first (this is typically the first thing filled in when
# logging a new purchase), then fill in the rest of the row's details.
$wsInventario.Hyperlinks.Add($wsInventario.Range("I12"), "https://tienda.bricogeek.com/servomotores/1320-mini-servo-feetech-3-5kg-ft1117m-fb-con-feedback.html")
$wsInventario.Range("G12").Value = "Fran"
$wsInventario.Range("E12").Value = "Servomotor realimentado"
$wsInventario.Range("C12").Value = "FT1117M-FB"
$wsInventario.Range("F12").Value = "Servomotores realimentados de 3.5kg" + [char]0xB7 + "cm"

# --- Compras: register the purchase made to pay for the servos above ---
$wsCompras.Range("B5").Value = 10
$wsCompras.Hyperlinks.Add($wsCompras.Range("B5"), "", "Inventario_cero_coste!B10", "", "10")
$wsCompras.Range("C5").Value = "Fran"
$wsCompras.Range("D5").Value = 25.96

# --- Update the selections left on each sheet ---
$wsPresupuesto.Range("C4").Select()
$wsCompras.Range("E5").Select()
$wsInventario.Activate()
$wsInventario.Range("B12").Select()
